# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the existing "sum" column (G) and filling in a 0 value
# for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting from G1 (bold font, thin border, centered)
# onto H1 so the new header looks consistent with the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
